$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.962.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").Value = "'3.757.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.71%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'595.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").Value = "'167.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.49%  "

$ws.Range("D7").Value = "'3.759.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.67%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.07%  "

$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.02%  "

$ws.Range("D11").Value = "'6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("D12").Value = "'0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.40%  "

$ws.Range("E13").Value = "  -7.66%  "

$ws.Range("D14").Value = "'36.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.09%  "

$ws.Range("D15").Value = "'4.386.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.76%  "

$ws.Range("D16").Value = "'3.752.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.03%  "

$ws.Range("D17").Value = "'67.806.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").Value = "'18.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").Value = "'7.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.63%  "

$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("D21").Value = "'10.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.01%  "

$ws.Range("D22").Value = "'465.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "

$ws.Range("D23").Value = "'0.708"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.88%  "

$ws.Range("D24").Value = "'83.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "

$ws.Range("E25").Value = "  -12.80%  "

$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.03%  "

$ws.Range("D27").Value = "'12.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "

$ws.Range("D28").Value = "'10.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").Value = "'3.903.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.77%  "

$ws.Range("E31").Value = "  -2.07%  "

$ws.Range("D32").Value = "'7.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.84%  "

$ws.Range("D33").Value = "'30.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("D34").Value = "'2.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.99%  "

$ws.Range("D35").Value = "'9.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.37%  "

$ws.Range("D36").Value = "'3.705.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.25%  "

$ws.Range("E37").Value = "  -3.45%  "

$ws.Range("D38").Value = "'3.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.24%  "

$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("D41").Value = "'5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.80%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'0.307"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.59%  "

$ws.Range("D45").Value = "'8.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.16%  "

$ws.Range("D46").Value = "'1.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.20%  "

$ws.Range("D47").Value = "'398.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.89%  "

$ws.Range("D48").Value = "'44.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.82%  "

$ws.Range("D49").Value = "'144.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("D50").Value = "'39.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("E51").Value = "  -3.67%  "
